# script changes done related to future rates
# Adds a new "futureRates" boolean flag column (J) to the info sheet:
#   J1 -> header label "futureRates" (styled like the other header cells)
#   J2 -> TRUE, displayed via a custom "TRUE"/"FALSE" number format

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, formatted the same way as the existing header row (I1).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("J1").Value = "futureRates"

# New boolean value cell with a custom TRUE/FALSE display format.
$j2 = $ws.Range("J2")
$j2.Value = $true
$j2.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# The header row's stale explicit height collapses back to the sheet default
# once the row is touched/resaved (matches the authored workbook).
$null = $ws.Rows.Item(1).AutoFit()

# Park the selection on the newly added cell, matching the authored workbook.
$null = $ws.Range("J2").Select()
